# LSTM_mv.xlsx update script
# Appends newly-realized actual values / newly generated predictions to the
# rolling daily model sheets, re-points the "trailing" placeholder row that
# only holds the latest known actual value, and moves the active sheet /
# selection to where the author left off.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "D1_USD" (sheet1) -----------------------------------------------
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("D1_USD")

# Row 90 already held only the realized C value; give it the matching
# prediction/date and the usual D/E formulas, then append the new rows.
$ws1.Range("A89:E89").Copy()
$ws1.Range("A90:E95").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A90").Value = 45282
$ws1.Range("B90").Value = 3.9279000000000002
$ws1.Range("D90").Formula = "=B90-C90"
$ws1.Range("E90").Formula = "=IF(D90<0,1,0)"

$ws1.Range("A91").Value = 45286
$ws1.Range("B91").Value = 3.9323109999999999
$ws1.Range("C91").Value = 3.8880050000000002
$ws1.Range("D91").Formula = "=B91-C91"
$ws1.Range("E91").Formula = "=IF(D91<0,1,0)"

$ws1.Range("A92").Value = 45287
$ws1.Range("B92").Value = 3.9164859999999999
$ws1.Range("C92").Value = 3.9465865999999998
$ws1.Range("D92").Formula = "=B92-C92"
$ws1.Range("E92").Formula = "=IF(D92<0,1,0)"

$ws1.Range("A93").Value = 45288
$ws1.Range("B93").Value = 3.886606
$ws1.Range("C93").Value = 3.9129738999999999
$ws1.Range("D93").Formula = "=B93-C93"
$ws1.Range("E93").Formula = "=IF(D93<0,1,0)"

$ws1.Range("A94").Value = 45289
$ws1.Range("B94").Value = 3.9155280000000001
$ws1.Range("C94").Value = 3.8855276000000001
$ws1.Range("D94").Formula = "=B94-C94"
$ws1.Range("E94").Formula = "=IF(D94<0,1,0)"

$ws1.Range("A95").Value = 45293
$ws1.Range("B95").Value = 3.9369000000000001
$ws1.Range("C95").Value = 3.9294790000000002
$ws1.Range("D95").Formula = "=B95-C95"
$ws1.Range("E95").Formula = "=IF(D95<0,1,0)"

# New trailing row: only the just-released actual is known so far.
$ws1.Range("C96").Value = 3.9566110000000001

# ---------------------------------------------------------------------------
# Sheet "D1_EUR" (sheet3) -------------------------------------------------
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("D1_EUR")

$ws3.Range("A363:E363").Copy()
$ws3.Range("A364:E369").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws3.Range("A364").Value = 45282
$ws3.Range("B364").Value = 4.3228
$ws3.Range("D364").Formula = "=B364-C364"
$ws3.Range("E364").Formula = "=IF(D364<0,1,0)"

$ws3.Range("A365").Value = 45286
$ws3.Range("B365").Value = 4.3332100000000002
$ws3.Range("C365").Value = 4.3285612999999996
$ws3.Range("D365").Formula = "=B365-C365"
$ws3.Range("E365").Formula = "=IF(D365<0,1,0)"

$ws3.Range("A366").Value = 45287
$ws3.Range("B366").Value = 4.3247400000000003
$ws3.Range("C366").Value = 4.3349209999999996
$ws3.Range("D366").Formula = "=B366-C366"
$ws3.Range("E366").Formula = "=IF(D366<0,1,0)"

$ws3.Range("A367").Value = 45288
$ws3.Range("B367").Value = 4.3173000000000004
$ws3.Range("C367").Value = 4.3282090000000002
$ws3.Range("D367").Formula = "=B367-C367"
$ws3.Range("E367").Formula = "=IF(D367<0,1,0)"

$ws3.Range("A368").Value = 45289
$ws3.Range("B368").Value = 4.3335499999999998
$ws3.Range("C368").Value = 4.3232400000000002
$ws3.Range("D368").Formula = "=B368-C368"
$ws3.Range("E368").Formula = "=IF(D368<0,1,0)"

$ws3.Range("A369").Value = 45293
$ws3.Range("B369").Value = 4.3479000000000001
$ws3.Range("C369").Value = 4.3370059999999997
$ws3.Range("D369").Formula = "=B369-C369"
$ws3.Range("E369").Formula = "=IF(D369<0,1,0)"

# New trailing row: only the just-released actual is known so far.
$ws3.Range("C370").Value = 4.3579197000000001

# ---------------------------------------------------------------------------
# Sheet "D5_EUR" (sheet5) -------------------------------------------------
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("D5_EUR")

# The day+5 prediction for these dates is now available, replacing the
# "Nan" placeholders that were waiting on it.
$ws5.Range("B55").Value = 4.32315
$ws5.Range("B56").Value = 4.34213
$ws5.Range("B57").Value = 4.3228
$ws5.Range("B58").Value = 4.3332100000000002
$ws5.Range("B59").Value = 4.3247400000000003

$ws5.Range("A59:C59").Copy()
$ws5.Range("A60:C64").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws5.Range("A60").Value = 45288
$ws5.Range("B60").Value = 4.3173000000000004
$ws5.Range("C60").Value = 4.2615129999999999

$ws5.Range("A61").Value = 45289
$ws5.Range("B61").Value = 4.3335499999999998
$ws5.Range("C61").Value = 4.2706203

$ws5.Range("A62").Value = 45293
$ws5.Range("B62").Value = 4.33988
$ws5.Range("C62").Value = 4.2816386

$ws5.Range("A63").Value = 45294
$ws5.Range("B63").Value = "Nan"
$ws5.Range("C63").Value = 4.2987776000000002

$ws5.Range("A64").Value = 45295
$ws5.Range("B64").Value = "Nan"
$ws5.Range("C64").Value = 4.3172812

# ---------------------------------------------------------------------------
# Sheet "D1_OIL" (sheet6) -------------------------------------------------
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("D1_OIL")

# Formulas are entered before the format is (re)copied from row 35 so the
# unstyled D/C columns keep their plain "General" look instead of inheriting
# the neighbouring currency-style formatting.
$ws6.Range("D36").Formula = "=B36-C36"
$ws6.Range("E36").Formula = "=D36/C36"
$ws6.Range("D37").Formula = "=B37-C37"
$ws6.Range("E37").Formula = "=D37/C37"
$ws6.Range("D38").Formula = "=B38-C38"
$ws6.Range("E38").Formula = "=D38/C38"
$ws6.Range("D39").Formula = "=B39-C39"
$ws6.Range("E39").Formula = "=D39/C39"

$ws6.Range("A35:E35").Copy()
$ws6.Range("A36:E39").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws6.Range("A35:B35").Copy()
$ws6.Range("A40:B41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws6.Range("A36").Value = 45282
$ws6.Range("B36").Value = 73.559997999999993
$ws6.Range("C36").Value = 71.102999999999994

$ws6.Range("A37").Value = 45286
$ws6.Range("B37").Value = 73.559997999999993
$ws6.Range("C37").Value = 67.181200000000004

$ws6.Range("A38").Value = 45287
$ws6.Range("B38").Value = 74.110000999999997
$ws6.Range("C38").Value = 67.802999999999997

$ws6.Range("A39").Value = 45288
$ws6.Range("B39").Value = 71.769997000000004
$ws6.Range("C39").Value = 68.910200000000003

$ws6.Range("A40").Value = 45289
$ws6.Range("B40").Value = 71.650002000000001

$ws6.Range("A41").Value = 45293
$ws6.Range("B41").Value = 70.379997000000003

# ---------------------------------------------------------------------------
# Restore each sheet's remembered selection. D1_USD is selected last so it
# ends up the active tab/sheet, matching the saved workbook.
# ---------------------------------------------------------------------------
$ws3.Range("C371").Select()
$ws5.Range("A62:B62").Select()
$ws6.Range("C40").Select()
$ws1.Range("C97").Select()
